$d = $word.ActiveDocument

# The target formula is the final equation in the document (the minimization
# objective that sums three cost terms: holding cost, separate-ordering
# cost, and joint-ordering cost). Locate it as the last OMath object.
$count = $d.OMaths.Count
$om = $d.OMaths.Item($count)

$xml = '<m:oMathPara><m:oMathParaPr><m:jc m:val="center" /></m:oMathParaPr><m:oMath><m:r><m:t>m</m:t></m:r><m:r><m:t>i</m:t></m:r><m:r><m:t>n</m:t></m:r><m:r><m:t>→</m:t></m:r><m:nary><m:naryPr><m:chr m:val="∑" /><m:limLoc m:val="undOvr" /><m:subHide m:val="0" /><m:supHide m:val="0" /></m:naryPr><m:sub><m:r><m:t>i</m:t></m:r><m:r><m:t>=</m:t></m:r><m:r><m:t>1</m:t></m:r></m:sub><m:sup><m:r><m:t>m</m:t></m:r><m:r><m:t>=</m:t></m:r><m:r><m:t>62</m:t></m:r></m:sup><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSub><m:e><m:r><m:t>q</m:t></m:r></m:e><m:sub><m:r><m:t>i</m:t></m:r></m:sub></m:sSub></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:nary><m:r><m:t>⋅</m:t></m:r><m:r><m:t>h</m:t></m:r><m:r><m:t>⋅</m:t></m:r><m:r><m:t>p</m:t></m:r><m:sSub><m:e><m:r><m:t>r</m:t></m:r></m:e><m:sub><m:r><m:t>i</m:t></m:r></m:sub></m:sSub><m:r><m:t>+</m:t></m:r><m:nary><m:naryPr><m:chr m:val="∑" /><m:limLoc m:val="undOvr" /><m:subHide m:val="0" /><m:supHide m:val="0" /></m:naryPr><m:sub><m:r><m:t>i</m:t></m:r><m:r><m:t>=</m:t></m:r><m:r><m:t>1</m:t></m:r></m:sub><m:sup><m:r><m:t>m</m:t></m:r><m:r><m:t>=</m:t></m:r><m:r><m:t>k</m:t></m:r></m:sup><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>d</m:t></m:r><m:sSub><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sub><m:r><m:t>i</m:t></m:r></m:sub></m:sSub></m:num><m:den><m:sSub><m:e><m:r><m:t>q</m:t></m:r></m:e><m:sub><m:r><m:t>i</m:t></m:r></m:sub></m:sSub></m:den></m:f></m:e></m:nary><m:r><m:t>⋅</m:t></m:r><m:sSup><m:e><m:r><m:t>c</m:t></m:r></m:e><m:sup><m:r><m:t>−</m:t></m:r><m:r><m:t>o</m:t></m:r><m:r><m:t>r</m:t></m:r></m:sup></m:sSup><m:r><m:t>+</m:t></m:r><m:nary><m:naryPr><m:chr m:val="∑" /><m:limLoc m:val="undOvr" /><m:subHide m:val="0" /><m:supHide m:val="0" /></m:naryPr><m:sub><m:r><m:t>i</m:t></m:r><m:r><m:t>=</m:t></m:r><m:r><m:t>1</m:t></m:r></m:sub><m:sup><m:r><m:t>m</m:t></m:r><m:r><m:t>=</m:t></m:r><m:r><m:t>62</m:t></m:r></m:sup><m:e><m:r><m:t>d</m:t></m:r></m:e></m:nary><m:sSub><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sub><m:r><m:t>i</m:t></m:r></m:sub></m:sSub><m:r><m:t>⋅</m:t></m:r><m:sSubSup><m:e><m:r><m:t>c</m:t></m:r></m:e><m:sub><m:r><m:t>i</m:t></m:r></m:sub><m:sup><m:r><m:t>o</m:t></m:r><m:r><m:t>r</m:t></m:r></m:sup></m:sSubSup></m:oMath></m:oMathPara>'

$om.Range.InsertXML($xml)

Write-Host ("Replaced OMath #" + $count)
